$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the _GoBack bookmark from its original position (start of
#    the "USE CASE PAYMENT" paragraph) - it will be re-added later at
#    the end of the new "Bank" paragraph.
# ------------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# ------------------------------------------------------------------
# 2. Paragraph 5 is currently empty with rPr (b, color=FF0000, sz24,
#    szCs24). Turn it into the "Secondary Actor:" line: change the
#    color to 0070C0 and add two runs: "Secondary" + " Actor:".
# ------------------------------------------------------------------
$pSecondary = $d.Paragraphs.Item(5)
$secStart = $pSecondary.Range.Start

$pSecondary.Range.InsertBefore("Secondary")
$r1 = $d.Range($secStart, $secStart + 9)
$r1.Font.Bold = $true
$r1.Font.Color = 0xC07000
$r1.Font.Size = 12
$r1.Font.SizeBi = 12

$r1.InsertAfter(" Actor:")
$r2 = $d.Range($secStart + 9, $secStart + 16)
$r2.Font.Bold = $true
$r2.Font.Color = 0xC07000
$r2.Font.Size = 12
$r2.Font.SizeBi = 12

# ------------------------------------------------------------------
# 3. Insert a new paragraph right after it for "Bank" (sz24/szCs24,
#    not bold, no color) and re-home the _GoBack bookmark at its end.
# ------------------------------------------------------------------
$pSecondary.Range.InsertParagraphAfter()
$pBank = $d.Paragraphs.Item(6)
$pBank.Range.InsertBefore("Bank")
$bankRun = $d.Range($pBank.Range.Start, $pBank.Range.Start + 4)
$bankRun.Font.Bold = $false
$bankRun.Font.Size = 12
$bankRun.Font.SizeBi = 12

$bmRange = $d.Range($pBank.Range.Start + 4, $pBank.Range.Start + 4)
$d.Bookmarks.Add("_GoBack", $bmRange)
